$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.685.52"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.52%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.887.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.20%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -1.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.42"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4828"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.86%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3786"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.49%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07325"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9173"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.86%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.41"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.70%  "

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.957.66"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.76%  "

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07680"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.46%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.462"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.587"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.73"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008788"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.723.05"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.50"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.112"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.145.01"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.80"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.72%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.904"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.76%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.44"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.79%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.34"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.85%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.102"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.82"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.896"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08925"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.151"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.45%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.219"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.95%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7569"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.619"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02031"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.29%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.542"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.70%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.17%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05245"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.16%  "

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5427"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.32%  "

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.968"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.41%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.941"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.30%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.82%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.319"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "109.39"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.99%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.61"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4769"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.37%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.626"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.22"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06054"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.69%  "
